$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Single-cell corrections ---
$ws.Range("Q55").Value = 0
$ws.Range("Q68").Value = 0
$ws.Range("Q69").Value = 0
$ws.Range("O1146").Value = 2
$ws.Range("R1148").Value = 0
$ws.Range("R1149").Value = 0

# --- New weekly rows 1150-1171 (columns A-Q; R left blank like existing new-data rows) ---
$newRows = @(
    @(45474, 120.8000030517578, 121.5500030517578, 117.3000030517578, 120.4499969482422, 120.4499969482422, 41483499, 2024, 7, 1, 0, 0, 0, 27, 0, 0, 0),
    @(45481, 120.5999984741211, 125.9000015258789, 119, 120.2900009155273, 120.2900009155273, 54081159, 2024, 7, 8, 0, 0, 0, 28, 0, 0, 0),
    @(45488, 120.5299987792969, 125.0899963378906, 119.6999969482422, 120.8899993896484, 120.8899993896484, 50871693, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 1),
    @(45495, 120, 122.4599990844727, 117, 119.370002746582, 119.370002746582, 37171924, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0),
    @(45502, 120.8499984741211, 128.1999969482422, 120.4000015258789, 126.2399978637695, 126.2399978637695, 58815147, 2024, 7, 29, 0, 0, 0, 31, 1, 0, 0),
    @(45509, 123.5999984741211, 124.9000015258789, 117.25, 118.9499969482422, 118.9499969482422, 55778142, 2024, 8, 5, 0, 0, 0, 32, 0, 0, 0),
    @(45516, 118.7900009155273, 119.5, 114, 116.0100021362305, 116.0100021362305, 22054298, 2024, 8, 12, 0, 0, 0, 33, 0, 0, 0),
    @(45523, 116.9499969482422, 121.5, 116.6900024414062, 120.5599975585938, 120.5599975585938, 20537565, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0),
    @(45530, 121, 121.4899978637695, 116.5, 117.7399978637695, 117.7399978637695, 26898805, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0),
    @(45537, 118.7900009155273, 119.3899993896484, 115.25, 115.75, 115.75, 23816095, 2024, 9, 2, 0, 0, 0, 36, 0, 0, 0),
    @(45544, 115, 115.5800018310547, 110.2200012207031, 112.3600006103516, 112.3600006103516, 30194921, 2024, 9, 9, 0, 0, 0, 37, 0, 0, 2),
    @(45551, 112.8899993896484, 115.25, 109.5100021362305, 110.3000030517578, 110.3000030517578, 21329408, 2024, 9, 16, 0, 0, 0, 38, 0, 0, 0),
    @(45558, 111.1500015258789, 113.5899963378906, 109.6600036621094, 110.7699966430664, 110.7699966430664, 26076097, 2024, 9, 23, 0, 0, 0, 39, 0, 0, 0),
    @(45565, 110.8199996948242, 111.6900024414062, 106.8000030517578, 108.7600021362305, 108.7600021362305, 17996004, 2024, 9, 30, 0, 0, 0, 40, 0, 0, 0),
    @(45572, 109, 109.4599990844727, 101.0699996948242, 105.5299987792969, 105.5299987792969, 19911763, 2024, 10, 7, 0, 0, 0, 41, 0, 0, 0),
    @(45579, 105.9100036621094, 106.3000030517578, 102.8199996948242, 104.8199996948242, 104.8199996948242, 14176542, 2024, 10, 14, 0, 0, 0, 42, 0, 0, 0),
    @(45586, 105.3000030517578, 106.1800003051758, 96, 96.55999755859375, 96.55999755859375, 33436521, 2024, 10, 21, 0, 0, 0, 43, 2, 0, 0),
    @(45593, 98, 110.7900009155273, 96.59999847412109, 110.1800003051758, 110.1800003051758, 36592826, 2024, 10, 28, 0, 0, 0, 44, 0, 0, 0),
    @(45600, 110.4400024414062, 115, 106.8600006103516, 111.6399993896484, 111.6399993896484, 30602569, 2024, 11, 4, 0, 0, 0, 45, 0, 0, 0),
    @(45607, 111.0299987792969, 115, 103.8000030517578, 104.3000030517578, 104.3000030517578, 31795371, 2024, 11, 11, 0, 0, 0, 46, 0, 0, 0),
    @(45614, 104.4000015258789, 107.5999984741211, 100.6999969482422, 103.0599975585938, 103.0599975585938, 22437060, 2024, 11, 18, 0, 0, 0, 47, 0, 0, 0),
    @(45621, 106.5899963378906, 114.6999969482422, 105.3399963378906, 110.5, 110.5, 41913231, 2024, 11, 25, 0, 0, 0, 48, 0, 0, 0)
)

$startRow = 1150
$r = $startRow
foreach ($row in $newRows) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $r = $r + 1
}

